# Daily attendance processing - 2026-01-04 14:58:27
# Reorders the "Recorded By" (column G) email lists for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"
    3  = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg"
    4  = "gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
    5  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    6  = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"
    7  = "NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    9  = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    11 = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    12 = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
    13 = "amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    15 = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    17 = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    19 = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    27 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    30 = "yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
